$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "version" + "+" (two adjacent bold/Georgia/18pt runs) collapse into
#    a single run "version+" - plain Find/Replace naturally coalesces
#    the two identically-formatted runs, which is what we want here.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("version+", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "version+", 2)

# ---------------------------------------------------------------------
# 2. "vector-size" -> "vector-length" inside a single run. The run
#    immediately before it happens to share identical run formatting,
#    and this engine (like Word) auto-coalesces an edited run into a
#    neighbor whose properties match. We don't want that merge here
#    (the target keeps the leading-space run separate), so we briefly
#    perturb the neighbor's formatting across the edit and restore it
#    afterwards, which prevents the accidental coalescing.
# ---------------------------------------------------------------------
$hit = $d.Content
$hit.Find.Execute("vector-size")
$start = $hit.Start
$end = $hit.End

$guard = $d.Range($start - 1, $start)
$originalSize = $guard.Font.Size
$guard.Font.Size = $originalSize + 1

$target = $d.Range($start, $end)
$target.Text = "vector-length"

$guard2 = $d.Range($start - 1, $start)
$guard2.Font.Size = $originalSize

# ---------------------------------------------------------------------
# 3. "use " + "mu" + "::{" (three adjacent same-format runs) collapse
#    into a single run "use mu::{" - again the natural coalescing
#    behaviour is exactly what the target wants.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("use mu::{", $true, $true, $false, $false, $false, `
                         $true, 1, $false, "use mu::{", 2)
